$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gof")
$ws.Range("D2").Value = 26061
$ws.Range("F2").Value = 26115
$ws.Range("G2").Value = 26259
$ws.Range("D3").Value = 26048
$ws.Range("F3").Value = 26134
$ws.Range("G3").Value = 26363

$ws = $wb.Worksheets.Item("facets")
$ws.Range("B2").Value = 754
$ws.Range("B3").Value = 746

$ws = $wb.Worksheets.Item("Estimates 0-1")
$ws.Range("B2").Value = 0.184
$ws.Range("E2").Value = 3.914
$ws.Range("B3").Value = -0.119
$ws.Range("D3").Value = -0.094
$ws.Range("E3").Value = 1.788
$ws.Range("B4").Value = 0.038
$ws.Range("D4").Value = 0.03
$ws.Range("E4").Value = 0.376
$ws.Range("B5").Value = -0.137
$ws.Range("C5").Value = 0.083
$ws.Range("D5").Value = -0.108
$ws.Range("E5").Value = 2.724
$ws.Range("B6").Value = 0.062
$ws.Range("C6").Value = 0.081
$ws.Range("D6").Value = 0.049
$ws.Range("E6").Value = 0.586
$ws.Range("B7").Value = 0.097
$ws.Range("C7").Value = 0.08
$ws.Range("D7").Value = 0.077
$ws.Range("E7").Value = 1.47
$ws.Range("B8").Value = 0.045
$ws.Range("D8").Value = 0.036
$ws.Range("E8").Value = 0.694
$ws.Range("B9").Value = -0.179
$ws.Range("C9").Value = 0.078
$ws.Range("D9").Value = -0.141
$ws.Range("E9").Value = 5.266
$ws.Range("B10").Value = -0.076
$ws.Range("C10").Value = 0.078
$ws.Range("D10").Value = -0.06
$ws.Range("E10").Value = 0.949
$ws.Range("B11").Value = -0.087
$ws.Range("C11").Value = 0.079
$ws.Range("D11").Value = -0.069
$ws.Range("E11").Value = 1.213
$ws.Range("B12").Value = 0.026
$ws.Range("C12").Value = 0.081
$ws.Range("D12").Value = 0.021
$ws.Range("E12").Value = 0.103
$ws.Range("B13").Value = -0.094
$ws.Range("C13").Value = 0.084
$ws.Range("D13").Value = -0.074
$ws.Range("E13").Value = 1.252
$ws.Range("B14").Value = 0.138
$ws.Range("C14").Value = 0.086
$ws.Range("D14").Value = 0.109
$ws.Range("E14").Value = 2.575
$ws.Range("B15").Value = -0.059
$ws.Range("C15").Value = 0.051
$ws.Range("D15").Value = -0.047
$ws.Range("E15").Value = 1.338
$ws.Range("B16").Value = -0.151
$ws.Range("C16").Value = 0.101
$ws.Range("D16").Value = -0.119
$ws.Range("E16").Value = 2.235
$ws.Range("B17").Value = -0.07
$ws.Range("D17").Value = -0.055
$ws.Range("E17").Value = 2.316
$ws.Range("B18").Value = -0.381
$ws.Range("C18").Value = 0.312
$ws.Range("D18").Value = -0.301
$ws.Range("E18").Value = 1.491

$ws = $wb.Worksheets.Item("Main effect 0-1")
$ws.Range("B2").Value = 0.067
$ws.Range("C2").Value = 0.053
$ws.Range("B3").Value = 0.088
$ws.Range("C3").Value = 0.069
